$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# Update test-case wording/typo fixes in the "String concatenate" and "String comparison" sections
$ws.Range("A4").Value = 'String concatenate:
strcat()'
$ws.Range("A5").Value = 'String concatenate:
strcat()'
$ws.Range("A6").Value = 'String concatenate:
strcat()'
$ws.Range("A7").Value = 'String concatenate:
strcat()'
$ws.Range("C7").Value = 'string1="Hello this is John Smith, Nice to meet you!"
string2="abc"'
$ws.Range("G7").Value = 'Causes: fgets function only reads BUFFER_SIZE -1 characters and will leave rest of characters in the input buffer.  Then when string2 call fgets function again, it read from input buffer directly. Also will has buffer overflow
Recommendation: Manually check the string1 length after uer input, if user exceed the limit prompt an error message and let the user input again'
$ws.Range("A8").Value = 'String concatenate:
strcat()'
$ws.Range("C8").Value = 'string1="12345" 
string2="Hi John Smith Nice to meet you!"'
$ws.Range("D8").Value = 'Output should be truncated: (12345Hi John Smith Nice to mee)
OR display an error message '
$ws.Range("G8").Value = 'Cause: fgets function truncated the string2 input base on the buffer size, however the combined string after strcat function leading to the buffer overflow
Recommendation: Manually check the string2 length after uer input, if user exceed the limit prompt an error message and let the user input again'
$ws.Range("A9").Value = 'String concatenate:
strcat()'
$ws.Range("D9").Value = 'Output should be truncated (Hello this is John Smith, Nice)
OR display an error message '
$ws.Range("G9").Value = 'Cause: In C, system strcat fucntion didn’t check the size of destination buffer(string1) automatically. It will simply append second string to the destination and past the end of buffer, leading to buffer overflow.This overflow can corrupt adjacent memory but might not immediately cause a crash.
Recommendation: Manually check the combined length of string is not greater than the buffer size before call strcat function, and if the combined string exceeds the limit, prompt an error message'
$ws.Range("A12").Value = 'String comparison: 
strcmp()'
$ws.Range("A13").Value = 'String comparison: 
strcmp()'
$ws.Range("A14").Value = 'String comparison: 
strcmp()'
$ws.Range("A15").Value = 'String comparison: 
strcmp()'
$ws.Range("A16").Value = 'String comparison: 
strcmp()'
$ws.Range("A17").Value = 'String comparison: 
strcmp()'
$ws.Range("G17").Value = 'Causes: fgets size limitation (detail can refer to G7)
Recommendation: Manually check the compare1 length after uer input, if user exceed the limit prompt an error message and let the user input again 
OR if the program will take the truncated input from compare1, a clean buffer function should be called after input of compare1 to allow input of compare2'
$ws.Range("A18").Value = 'String comparison: 
strcmp()'
$ws.Range("D18").Value = 'compare2 should be truncated to BUFFER_SIZE-1, and the output should be compare1 and compare2 is same ("Hello this is John Smith, Nice" string is equal to "Hello this is John Smith, Nice")
Can''t input the next compare1 value'
$ws.Range("G18").Value = 'Causes: fgets size limitation (detail can refer to G7)
Recommendation: Manually check the compare2 length after uer input, if user exceed the limit prompt an error message and let the user input again
OR if the program will take the truncated input from compare2, a clean buffer function should be called after input of compare2 to allow input of next compare1'

# Restore view / selection state
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("H15").Select()
